$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns E, F, G, H for rows 2 through 10 with new computed values.
$values = @{
    2  = @(1,                 0.2300463546529723, 1.47312266573296,   4.175686274509804)
    3  = @(1.333333333333333, 0.9035628306878307, 3.804974794238683,  9.017962962962963)
    4  = @(1.132075471698113, 0.6101235212638515, 2.972091344713986,  7.520754716981132)
    5  = @(1,                 0.3143800329074719, 3.824024366047232,  8.706829268292683)
    6  = @(1.428571428571429, 1.133932985098801,  4.242640954000648,  9.810204081632653)
    7  = @(1.209302325581395, 0.8171995893318567, 4.336273301956441,  9.983488372093023)
    8  = @(1,                 0.1253704264191552, 0.5992779459241324, 2.119322033898305)
    9  = @(1.2,                0.5421815972222221, 2.041969113756613,  5.497833333333333)
    10 = @(1.05,               0.3355795469576719, 1.598878819444445,  4.365166666666667)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("E$row").Value = $rowVals[0]
    $ws.Range("F$row").Value = $rowVals[1]
    $ws.Range("G$row").Value = $rowVals[2]
    $ws.Range("H$row").Value = $rowVals[3]
}
